# DPLKKPS133_RegisTamPeserta.xlsx
#
# The commit regenerates the "No. Rekening" (account number) used by the
# DPLKKPS133-002 ("Lanjutkan ke Verifikasi") test case, and leaves the
# workbook with that sheet active/selected (matching where the author was
# last working), instead of DPLKKPS133-001.
#
# Concretely:
#   - Sheet "DPLKKPS133-002", cell P2 (NOMOR_REKENING): 790854187 -> 790850353
#   - Sheet "DPLKKPS133-002", cell F2 (PREPARATION) contains the same
#     account number embedded in a multi-line text block; update it too.
#   - The active sheet/tab becomes DPLKKPS133-002 (was DPLKKPS133-001),
#     with the selection moved to Q2 on that sheet.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("DPLKKPS133-002")

# Update the numeric account-number cell.
$ws2.Range("P2").Value = 790850353

# Update the same account number embedded inside the PREPARATION text block.
$prep = $ws2.Range("F2").Value2
$ws2.Range("F2").Value = $prep.Replace("790854187", "790850353")

# Make DPLKKPS133-002 the active sheet/tab, with Q2 selected.
$ws2.Activate()
$ws2.Range("Q2").Select()
